$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1725663716814159
$ws.Range("C2").Value = 0.6238938053097345
$ws.Range("P2").Value = 0.1283185840707965
$ws.Range("S2").Value = 0.0752212389380531
$ws.Range("B3").Value = 0.0136986301369863
$ws.Range("C3").Value = 0.0273972602739726
$ws.Range("J3").Value = 0.04794520547945205
$ws.Range("P3").Value = 0.6575342465753424
$ws.Range("S3").Value = 0.2534246575342466
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("P4").Value = 0.7446808510638298
$ws.Range("S4").Value = 0.2127659574468085
$ws.Range("B6").Value = 0.03755868544600939
$ws.Range("D6").Value = 0.0187793427230047
$ws.Range("F6").Value = 0.03755868544600939
$ws.Range("J6").Value = 0.244131455399061
$ws.Range("O6").Value = 0.02816901408450704
$ws.Range("Q6").Value = 0.1690140845070423
$ws.Range("R6").Value = 0.07511737089201878
$ws.Range("S6").Value = 0.3896713615023474
$ws.Range("B7").Value = 0.07657657657657657
$ws.Range("D7").Value = 0.02252252252252252
$ws.Range("F7").Value = 0.04954954954954955
$ws.Range("J7").Value = 0.1126126126126126
$ws.Range("O7").Value = 0.02702702702702703
$ws.Range("Q7").Value = 0.2027027027027027
$ws.Range("R7").Value = 0.1171171171171171
$ws.Range("S7").Value = 0.3918918918918919
$ws.Range("B8").Value = 0.09178743961352658
$ws.Range("D8").Value = 0.01690821256038647
$ws.Range("F8").Value = 0.05314009661835749
$ws.Range("J8").Value = 0.1352657004830918
$ws.Range("O8").Value = 0.03381642512077294
$ws.Range("Q8").Value = 0.1908212560386473
$ws.Range("R8").Value = 0.1400966183574879
$ws.Range("S8").Value = 0.3381642512077295
$ws.Range("B9").Value = 0.0855614973262032
$ws.Range("D9").Value = 0.0267379679144385
$ws.Range("F9").Value = 0.06951871657754011
$ws.Range("J9").Value = 0.1283422459893048
$ws.Range("O9").Value = 0.0106951871657754
$ws.Range("Q9").Value = 0.1711229946524064
$ws.Range("R9").Value = 0.106951871657754
$ws.Range("S9").Value = 0.4010695187165775
$ws.Range("B10").Value = 0.08347107438016529
$ws.Range("D10").Value = 0.02148760330578512
$ws.Range("E10").Value = 0.0008264462809917355
$ws.Range("F10").Value = 0.07933884297520662
$ws.Range("J10").Value = 0.1049586776859504
$ws.Range("O10").Value = 0.0256198347107438
$ws.Range("Q10").Value = 0.2264462809917355
$ws.Range("R10").Value = 0.09090909090909091
$ws.Range("S10").Value = 0.3669421487603306
$ws.Range("G11").Value = 0.1517857142857143
$ws.Range("J11").Value = 0.06845238095238096
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.5744047619047619
$ws.Range("S11").Value = 0.01785714285714286
$ws.Range("G12").Value = 0.7577319587628866
$ws.Range("J12").Value = 0.1804123711340206
$ws.Range("K12").Value = 0.0154639175257732
$ws.Range("L12").Value = 0.01030927835051546
$ws.Range("S12").Value = 0.03608247422680412
$ws.Range("F13").Value = 0.02173913043478261
$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2826086956521739
$ws.Range("F15").Value = 0.01652892561983471
$ws.Range("H15").Value = 0.140495867768595
$ws.Range("I15").Value = 0.06611570247933884
$ws.Range("J15").Value = 0.3305785123966942
$ws.Range("K15").Value = 0.06198347107438017
$ws.Range("M15").Value = 0.01239669421487603
$ws.Range("O15").Value = 0.04958677685950413
$ws.Range("S15").Value = 0.3223140495867768
$ws.Range("H16").Value = 0.2384105960264901
$ws.Range("I16").Value = 0.03973509933774835
$ws.Range("J16").Value = 0.4437086092715232
$ws.Range("K16").Value = 0.08609271523178808
$ws.Range("M16").Value = 0.02649006622516556
$ws.Range("N16").Value = 0.006622516556291391
$ws.Range("O16").Value = 0.06622516556291391
$ws.Range("S16").Value = 0.09271523178807947
$ws.Range("F17").Value = 0.01505376344086022
$ws.Range("H17").Value = 0.1526881720430107
$ws.Range("I17").Value = 0.08387096774193549
$ws.Range("J17").Value = 0.4258064516129032
$ws.Range("K17").Value = 0.1247311827956989
$ws.Range("M17").Value = 0.02365591397849462
$ws.Range("N17").Value = 0.002150537634408602
$ws.Range("O17").Value = 0.06881720430107527
$ws.Range("S17").Value = 0.1032258064516129
$ws.Range("F18").Value = 0.02620087336244541
$ws.Range("H18").Value = 0.1965065502183406
$ws.Range("I18").Value = 0.1004366812227074
$ws.Range("J18").Value = 0.3406113537117904
$ws.Range("K18").Value = 0.1135371179039301
$ws.Range("M18").Value = 0.03493449781659388
$ws.Range("O18").Value = 0.07860262008733625
$ws.Range("S18").Value = 0.1091703056768559
$ws.Range("F19").Value = 0.01516427969671441
$ws.Range("H19").Value = 0.1946082561078349
$ws.Range("I19").Value = 0.08930075821398484
$ws.Range("J19").Value = 0.367312552653749
$ws.Range("K19").Value = 0.1297388374052232
$ws.Range("M19").Value = 0.01769165964616681
$ws.Range("N19").Value = 0.0008424599831508003
$ws.Range("O19").Value = 0.07245155855096883
$ws.Range("S19").Value = 0.1128896377422072
